$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing shared strings (rename "... Mux (CD4051BE)" -> "... Mux Out (CD4051BE)") ---
# These values live in F12 and F14 respectively.
$ws.Range("F12").Value = "PIP Analog Mux Out (CD4051BE)"
$ws.Range("F14").Value = "MCP Analog Mux Out (CD4051BE)"

# --- Clear the "Toggle Switch / Power Regulator" entry from row 5 (E5:F5) ---
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""

# --- Clear the "Blue Rail / Negative Power Rail" entry from row 6 (E6:F6) ---
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""

# --- Add a new "RX" label to L14 ---
$ws.Range("L14").Value = "RX"

# --- Move the "Toggle Switch +5V" / "Power Regulator Circuit" entry into row 17 (E17:F17) ---
$ws.Range("E17").Value = "Toggle Switch +5V"
$ws.Range("F17").Value = "Power Regulator Circuit"
# F17 also picks up F5's original cell format (fill/border) - copy formats only
$ws.Range("F5").Copy()
$ws.Range("F17").PasteSpecial(-4122)

# --- Move the "Blue Rail" / "Negative Power Rail" entry into row 18 (E18:F18) ---
$ws.Range("E18").Value = "Blue Rail"
$ws.Range("F18").Value = "Negative Power Rail"

# --- New Bluetooth RS232 RX/TX entries in rows 19-20 (columns A:B) ---
$ws.Range("A19").Value = "Bluetooth RS232"
$ws.Range("B19").Value = "TX"
$ws.Range("A20").Value = "Bluetooth RS232"
$ws.Range("B20").Value = "RX"

# --- Update the active selection to F21 ---
$ws.Range("F21").Select()
